$wb = $excel.ActiveWorkbook

# --- Work on the original (data) sheet first ---
$ws1 = $wb.Worksheets.Item(1)

# Insert a new column before column A for the "ID" field, shifting existing
# headers (Rok, Miesiac, MlodyWest, Hubert., diff, <-- West - Hubert) right.
$ws1.Columns.Item(1).Insert()
$ws1.Range("A1").Value = "ID"

# Insert a new column before "MlodyWest" (now column D after the shift above)
# for the "Dzien" field.
$ws1.Columns.Item(4).Insert()
$ws1.Range("D1").Value = "Dzien"

# Fill in the new data row (row 2) with one day's worth of records.
$ws1.Range("A2").Value = 2138
$ws1.Range("B2").Value = 2024
$ws1.Range("C2").Value = 9
$ws1.Range("D2").Value = 30
$ws1.Range("E2").Value = 154720
$ws1.Range("F2").Value = 148783
$ws1.Range("G2").Value = 5937

# Rename the data sheet and move the selection off the edited area.
$ws1.Name = "Liczby"
$ws1.Range("G8").Select()

# --- Add a second sheet to host a chart ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Chart"
$ws2.Range("A2").Select()
